# Natmi following Dr Hou advice
# Rewrite the LR-pair result matrix (rows 2-21) with updated stats
# using 3 replicate cells (was 1) and the full 5x4 sending/target cluster grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20
$data[0,0] = "ECs"
$data[0,1] = "Sema4a"
$data[0,2] = "Plxnb1"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 8.242326666666667
$data[0,7] = 24.72698
$data[0,8] = 0.1455255597722841
$data[0,9] = 0.1455255597722841
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.9392856666666667
$data[0,13] = 2.817857
$data[0,14] = 0.193119603916127
$data[0,15] = 0.193119603916127
$data[0,16] = 7.741899297984445
$data[0,17] = 69.67709368186
$data[0,18] = 0.02810383846289617
$data[0,19] = 0.02810383846289617

$data[1,0] = "ECs"
$data[1,1] = "Sema4a"
$data[1,2] = "Plxnb1"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 8.242326666666667
$data[1,7] = 24.72698
$data[1,8] = 0.1455255597722841
$data[1,9] = 0.1455255597722841
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 2.006819333333333
$data[1,13] = 6.020458
$data[1,14] = 0.4126073339966074
$data[1,15] = 0.4126073339966074
$data[1,16] = 16.54086050631555
$data[1,17] = 148.86774455684
$data[1,18] = 0.06004491324600606
$data[1,19] = 0.06004491324600607

$data[2,0] = "ECs"
$data[2,1] = "Sema4a"
$data[2,2] = "Plxnb1"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 8.242326666666667
$data[2,7] = 24.72698
$data[2,8] = 0.1455255597722841
$data[2,9] = 0.1455255597722841
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.01023466666666667
$data[2,13] = 0.030704
$data[2,14] = 0.002104274389594917
$data[2,15] = 0.002104274389594917
$data[2,16] = 0.08435746599111112
$data[2,17] = 0.75921719392
$data[2,18] = 0.0003062257084602816
$data[2,19] = 0.0003062257084602817

$data[3,0] = "ECs"
$data[3,1] = "Sema4a"
$data[3,2] = "Plxnb1"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 8.242326666666667
$data[3,7] = 24.72698
$data[3,8] = 0.1455255597722841
$data[3,9] = 0.1455255597722841
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 1.907411333333333
$data[3,13] = 5.722234
$data[3,14] = 0.3921687876976707
$data[3,15] = 0.3921687876976707
$data[3,16] = 15.72150729703556
$data[3,17] = 141.49356567332
$data[3,18] = 0.05707058235492155
$data[3,19] = 0.05707058235492157

$data[4,0] = "FAPs"
$data[4,1] = "Sema4a"
$data[4,2] = "Plxnb1"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 16.22192133333333
$data[4,7] = 48.665764
$data[4,8] = 0.2864123539488392
$data[4,9] = 0.2864123539488393
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.9392856666666667
$data[4,13] = 2.817857
$data[4,14] = 0.193119603916127
$data[4,15] = 0.193119603916127
$data[4,16] = 15.23701819419422
$data[4,17] = 137.133163747748
$data[4,18] = 0.05531184035128542
$data[4,19] = 0.05531184035128542

$data[5,0] = "FAPs"
$data[5,1] = "Sema4a"
$data[5,2] = "Plxnb1"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 16.22192133333333
$data[5,7] = 48.665764
$data[5,8] = 0.2864123539488392
$data[5,9] = 0.2864123539488393
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 2.006819333333333
$data[5,13] = 6.020458
$data[5,14] = 0.4126073339966074
$data[5,15] = 0.4126073339966074
$data[5,16] = 32.55446535554577
$data[5,17] = 292.990188199912
$data[5,18] = 0.1181758377865232
$data[5,19] = 0.1181758377865233

$data[6,0] = "FAPs"
$data[6,1] = "Sema4a"
$data[6,2] = "Plxnb1"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 16.22192133333333
$data[6,7] = 48.665764
$data[6,8] = 0.2864123539488392
$data[6,9] = 0.2864123539488393
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.01023466666666667
$data[6,13] = 0.030704
$data[6,14] = 0.002104274389594917
$data[6,15] = 0.002104274389594917
$data[6,16] = 0.1660259575395555
$data[6,17] = 1.494233617856
$data[6,18] = 0.0006026901812781369
$data[6,19] = 0.000602690181278137

$data[7,0] = "FAPs"
$data[7,1] = "Sema4a"
$data[7,2] = "Plxnb1"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 16.22192133333333
$data[7,7] = 48.665764
$data[7,8] = 0.2864123539488392
$data[7,9] = 0.2864123539488393
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 1.907411333333333
$data[7,13] = 5.722234
$data[7,14] = 0.3921687876976707
$data[7,15] = 0.3921687876976707
$data[7,16] = 30.94187659964177
$data[7,17] = 278.476889396776
$data[7,18] = 0.1123219856297524
$data[7,19] = 0.1123219856297525

$data[8,0] = "M1"
$data[8,1] = "Sema4a"
$data[8,2] = "Plxnb1"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 10.05778
$data[8,7] = 30.17334
$data[8,8] = 0.1775789924082702
$data[8,9] = 0.1775789924082702
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.9392856666666667
$data[8,13] = 2.817857
$data[8,14] = 0.193119603916127
$data[8,15] = 0.193119603916127
$data[8,16] = 9.447128592486667
$data[8,17] = 85.02415733238
$data[8,18] = 0.03429398467771007
$data[8,19] = 0.03429398467771008

$data[9,0] = "M1"
$data[9,1] = "Sema4a"
$data[9,2] = "Plxnb1"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 10.05778
$data[9,7] = 30.17334
$data[9,8] = 0.1775789924082702
$data[9,9] = 0.1775789924082702
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 2.006819333333333
$data[9,13] = 6.020458
$data[9,14] = 0.4126073339966074
$data[9,15] = 0.4126073339966074
$data[9,16] = 20.18414735441333
$data[9,17] = 181.65732618972
$data[9,18] = 0.07327039463138015
$data[9,19] = 0.07327039463138016

$data[10,0] = "M1"
$data[10,1] = "Sema4a"
$data[10,2] = "Plxnb1"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 10.05778
$data[10,7] = 30.17334
$data[10,8] = 0.1775789924082702
$data[10,9] = 0.1775789924082702
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.01023466666666667
$data[10,13] = 0.030704
$data[10,14] = 0.002104274389594917
$data[10,15] = 0.002104274389594917
$data[10,16] = 0.1029380257066667
$data[10,17] = 0.9264422313599999
$data[10,18] = 0.0003736749258547931
$data[10,19] = 0.0003736749258547932

$data[11,0] = "M1"
$data[11,1] = "Sema4a"
$data[11,2] = "Plxnb1"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 10.05778
$data[11,7] = 30.17334
$data[11,8] = 0.1775789924082702
$data[11,9] = 0.1775789924082702
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 1.907411333333333
$data[11,13] = 5.722234
$data[11,14] = 0.3921687876976707
$data[11,15] = 0.3921687876976707
$data[11,16] = 19.18432356017333
$data[11,17] = 172.65891204156
$data[11,18] = 0.06964093817332519
$data[11,19] = 0.0696409381733252

$data[12,0] = "M2"
$data[12,1] = "Sema4a"
$data[12,2] = "Plxnb1"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 19.67600733333333
$data[12,7] = 59.028022
$data[12,8] = 0.3473972941216719
$data[12,9] = 0.347397294121672
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 0.9392856666666667
$data[12,13] = 2.817857
$data[12,14] = 0.193119603916127
$data[12,15] = 0.193119603916127
$data[12,16] = 18.48139166542822
$data[12,17] = 166.332524988854
$data[12,18] = 0.06708922784231157
$data[12,19] = 0.06708922784231157

$data[13,0] = "M2"
$data[13,1] = "Sema4a"
$data[13,2] = "Plxnb1"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 19.67600733333333
$data[13,7] = 59.028022
$data[13,8] = 0.3473972941216719
$data[13,9] = 0.347397294121672
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 2.006819333333333
$data[13,13] = 6.020458
$data[13,14] = 0.4126073339966074
$data[13,15] = 0.4126073339966074
$data[13,16] = 39.48619191934178
$data[13,17] = 355.375727274076
$data[13,18] = 0.1433386713651783
$data[13,19] = 0.1433386713651784

$data[14,0] = "M2"
$data[14,1] = "Sema4a"
$data[14,2] = "Plxnb1"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 19.67600733333333
$data[14,7] = 59.028022
$data[14,8] = 0.3473972941216719
$data[14,9] = 0.347397294121672
$data[14,10] = 1
$data[14,11] = 0.3333333333333333
$data[14,12] = 0.01023466666666667
$data[14,13] = 0.030704
$data[14,14] = 0.002104274389594917
$data[14,15] = 0.002104274389594917
$data[14,16] = 0.2013773763875556
$data[14,17] = 1.812396387488
$data[14,18] = 0.0007310192290348069
$data[14,19] = 0.000731019229034807

$data[15,0] = "M2"
$data[15,1] = "Sema4a"
$data[15,2] = "Plxnb1"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 19.67600733333333
$data[15,7] = 59.028022
$data[15,8] = 0.3473972941216719
$data[15,9] = 0.347397294121672
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 1.907411333333333
$data[15,13] = 5.722234
$data[15,14] = 0.3921687876976707
$data[15,15] = 0.3921687876976707
$data[15,16] = 37.53023938234978
$data[15,17] = 337.772154441148
$data[15,18] = 0.1362383756851472
$data[15,19] = 0.1362383756851472

$data[16,0] = "sCs"
$data[16,1] = "Sema4a"
$data[16,2] = "Plxnb1"
$data[16,3] = "ECs"
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 2.440308333333334
$data[16,7] = 7.320925000000001
$data[16,8] = 0.04308579974893452
$data[16,9] = 0.04308579974893453
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 0.9392856666666667
$data[16,13] = 2.817857
$data[16,14] = 0.193119603916127
$data[16,15] = 0.193119603916127
$data[16,16] = 2.292146639747223
$data[16,17] = 20.629319757725
$data[16,18] = 0.0083207125819238
$data[16,19] = 0.0083207125819238

$data[17,0] = "sCs"
$data[17,1] = "Sema4a"
$data[17,2] = "Plxnb1"
$data[17,3] = "FAPs"
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 2.440308333333334
$data[17,7] = 7.320925000000001
$data[17,8] = 0.04308579974893452
$data[17,9] = 0.04308579974893453
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 2.006819333333333
$data[17,13] = 6.020458
$data[17,14] = 0.4126073339966074
$data[17,15] = 0.4126073339966074
$data[17,16] = 4.897257942627778
$data[17,17] = 44.07532148365
$data[17,18] = 0.01777751696751957
$data[17,19] = 0.01777751696751957

$data[18,0] = "sCs"
$data[18,1] = "Sema4a"
$data[18,2] = "Plxnb1"
$data[18,3] = "M2"
$data[18,4] = 3
$data[18,5] = 1
$data[18,6] = 2.440308333333334
$data[18,7] = 7.320925000000001
$data[18,8] = 0.04308579974893452
$data[18,9] = 0.04308579974893453
$data[18,10] = 1
$data[18,11] = 0.3333333333333333
$data[18,12] = 0.01023466666666667
$data[18,13] = 0.030704
$data[18,14] = 0.002104274389594917
$data[18,15] = 0.002104274389594917
$data[18,16] = 0.02497574235555556
$data[18,17] = 0.2247816812
$data[18,18] = 0.00009066434496689801
$data[18,19] = 0.00009066434496689802

$data[19,0] = "sCs"
$data[19,1] = "Sema4a"
$data[19,2] = "Plxnb1"
$data[19,3] = "sCs"
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 2.440308333333334
$data[19,7] = 7.320925000000001
$data[19,8] = 0.04308579974893452
$data[19,9] = 0.04308579974893453
$data[19,10] = 3
$data[19,11] = 1
$data[19,12] = 1.907411333333333
$data[19,13] = 5.722234
$data[19,14] = 0.3921687876976707
$data[19,15] = 0.3921687876976707
$data[19,16] = 4.654671771827779
$data[19,17] = 41.89204594645
$data[19,18] = 0.01689690585452425
$data[19,19] = 0.01689690585452426

$ws.Range("A2:T21").Value = $data
